# Rename the "TreatmentSupesedeRules" sheet to "Treatment Supersede Rules"
# (fixes the typo / adds proper spacing in the sheet's tab name).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TreatmentSupesedeRules")
$ws.Name = "Treatment Supersede Rules"
